$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.690.72"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "1.871.67"
$ws.Range("E3").Value = "  +2.30%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.38"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4617"
$ws.Range("E7").Value = "  -0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3888"
$ws.Range("E8").Value = "  +0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07878"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9762"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.00"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "1.839.08"
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.010"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.706"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06968"
$ws.Range("E15").Value = "  +1.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.37"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001003"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.85"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").Value = "28.676.18"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.276"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.10"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "2.078.99"
$ws.Range("E25").Value = "  -2.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.75"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.30"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.866"
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.989"
$ws.Range("E29").Value = "  +1.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.46"
$ws.Range("E30").Value = "  +1.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09328"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9189"
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.285"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.335"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.322"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05800"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02079"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.691"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5633"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1786"
$ws.Range("E41").Value = "  +1.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.780"
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07227"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.71"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5290"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.159"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.123"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.839"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.98"
$ws.Range("E49").Value = "  +0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.410"
$ws.Range("E50").Value = "  +3.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.25%  "
